# Update LeagueData.xlsx for the Feb 2026 offsets refresh (discobisco offsets).
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "BasePointers"
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("BasePointers")

# "History" and "Jersey" pointers no longer exist - drop both rows (they are
# contiguous: row 4 = History, row 5 = Jersey).
$ws1.Rows("4:5").Delete()

# "career_stats" pointer is gone too. After the deletion above it now lives
# at row 10 (was row 12).
$ws1.Rows("10").Delete()

# Refresh the remaining pointer offsets (2k26_offsets.json column). The
# leading "'" keeps these long numbers stored as literal text (matching the
# existing cells) instead of being coerced to numeric values.
$ws1.Range("F3").Value = "'5348811252"   # HallOfFame
$ws1.Range("F4").Value = "'5347512576"   # NBAHistory
$ws1.Range("F5").Value = "'131015096"    # Player
$ws1.Range("F6").Value = "'130992024"    # Stadium
# Staff (row 7) offset is unchanged (130991496) - nothing to update.
$ws1.Range("F8").Value = "'131015696"    # Team
$ws1.Range("F9").Value = "'5336537040"   # TeamHistory

# ---------------------------------------------------------------------------
# Sheet "GameInfo"
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("GameInfo")

# career_statsSize, coachSize, hall_of_fameSize, historySize and jerseySize
# fields are gone. Rows 2 and 3 are contiguous (career_statsSize/coachSize).
# After removing them, hall_of_fameSize/historySize/jerseySize shift up to
# rows 3:5.
$ws2.Rows("2:3").Delete()
$ws2.Rows("3:5").Delete()

# Only the build date (2k26_offsets.json version row) actually changed.
$ws2.Range("F7").Value = "'2026-02-22"
